$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: Title paragraph "Task Breakdown Document" -> append a new
# run containing " - Phase 2" (en dash), as its own <w:r>.
# ---------------------------------------------------------------------
$titleRange = $d.Content
$titleRange.Find.Execute("Task Breakdown Document", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$titleEnd = $titleRange.End

$appendRange = $d.Range($titleEnd, $titleEnd)
$appendRange.InsertAfter(" – Phase 2")

# Force the newly-typed text into its own run (distinct from the
# existing "Task Breakdown Document" run) the same way Word does when
# formatting is (re)applied to a sub-range: toggle Bold on then off.
$newTextRange = $d.Range($titleEnd, $titleEnd + 10)
$newTextRange.Font.Bold = $true
$newTextRange.Font.Bold = $false

# ---------------------------------------------------------------------
# Change 2: "Created add system for cart" table cell -> split into
# three runs ("Created " / "add system" / " for cart") with a
# w:proofErr gramStart/gramEnd pair bracketing "add system", matching
# the markers Word's grammar checker leaves behind.
# ---------------------------------------------------------------------
$cellRange = $d.Content
$cellRange.Find.Execute("Created add system for cart", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cellStart = $cellRange.Start
$cellEnd = $cellRange.End

$paraXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5FB69BEE" w14:textId="4A977650" w:rsidR="00A74731" w:rsidRDefault="00A74731" w:rsidP="00A74731"><w:pPr><w:contextualSpacing/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">Created </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>add system</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> for cart</w:t></w:r></w:p>'

$cellWhole = $d.Range($cellStart, $cellEnd)
$cellWhole.InsertXML($paraXml) | Out-Null

Write-Output "edit complete"
